# "fix sel card bug"
#
# The "Rarity" column header (C1) is clarified to spell out what the
# numeric values mean, and the leftover cell selection (from scrolling
# past the data while testing) is corrected back to the row right below
# the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Rarity (0 = Common, 1 = Rare, 2 =VeryRare, 3 = Epic)"

$ws.Range("E10").Select()
